$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Swap the content of rows 7 and 8 for the columns that differ between them.
# (Columns A, B, D, E, F, G, H, Q, R, Y, AA)

# --- Save row 7's original values ---
$A7 = $ws.Range("A7").Value2
$B7 = $ws.Range("B7").Value2
$D7 = $ws.Range("D7").Value2
$E7 = $ws.Range("E7").Value2
$F7 = $ws.Range("F7").Value2
$G7 = $ws.Range("G7").Value2
$H7 = $ws.Range("H7").Value2
$Q7 = $ws.Range("Q7").Value2
$R7 = $ws.Range("R7").Value2
$Y7 = $ws.Range("Y7").Value2
$AA7 = $ws.Range("AA7").Value2

# --- Save row 8's original values ---
$A8 = $ws.Range("A8").Value2
$B8 = $ws.Range("B8").Value2
$D8 = $ws.Range("D8").Value2
$E8 = $ws.Range("E8").Value2
$F8 = $ws.Range("F8").Value2
$G8 = $ws.Range("G8").Value2
$H8 = $ws.Range("H8").Value2
$Q8 = $ws.Range("Q8").Value2
$R8 = $ws.Range("R8").Value2
$Y8 = $ws.Range("Y8").Value2
$AA8 = $ws.Range("AA8").Value2

# --- Write row 8's original values into row 7 ---
$ws.Range("A7").Value2 = $A8
$ws.Range("B7").Value2 = $B8
$ws.Range("D7").Value2 = $D8
$ws.Range("E7").Value2 = $E8
$ws.Range("F7").Value2 = $F8
$ws.Range("G7").Value2 = $G8
$ws.Range("H7").Value2 = $H8
$ws.Range("Q7").Value2 = $Q8
$ws.Range("R7").Value2 = $R8

# Y/AA hold date-like text ("2023-09-07"); a leading apostrophe forces Excel
# to keep them as literal text instead of auto-converting to a date serial,
# and ClearFormats afterwards drops the quote-prefix cell style so the cell
# ends up with the plain default style, matching the original text cells.
$ws.Range("Y7").Value2 = "'" + $Y8
$ws.Range("Y7").ClearFormats()
$ws.Range("AA7").Value2 = "'" + $AA8
$ws.Range("AA7").ClearFormats()

# --- Write row 7's original values into row 8 ---
$ws.Range("A8").Value2 = $A7
$ws.Range("B8").Value2 = $B7
$ws.Range("D8").Value2 = $D7
$ws.Range("E8").Value2 = $E7
$ws.Range("F8").Value2 = $F7
$ws.Range("G8").Value2 = $G7
$ws.Range("H8").Value2 = $H7
$ws.Range("Q8").Value2 = $Q7
$ws.Range("R8").Value2 = $R7

$ws.Range("Y8").Value2 = "'" + $Y7
$ws.Range("Y8").ClearFormats()
$ws.Range("AA8").Value2 = "'" + $AA7
$ws.Range("AA8").ClearFormats()
